$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 34) down to the new row 35
$ws.Range("A34:G34").Copy()
$ws.Range("A35:G35").PasteSpecial(-4122)

$ws.Range("A35").Value = "4/18/2020"
$ws.Range("B35").Value = -754
$ws.Range("C35").ClearContents()
$ws.Range("D35").Value = -112
$ws.Range("E35").ClearContents()
$ws.Range("F35").Value = 507
$ws.Range("G35").Value = 1384

# Match the cursor position left behind by the editor (one row below the new data)
[void]$ws.Range("F36").Select()
